$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MakeAppointmentIntend")
$ws.Rows("5:6").Delete()
